$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.443.55'
$ws.Range('E2').Value = '  +0.09%  '
$ws.Range('D3').Value = '1.939.44'
$ws.Range('E3').Value = '  -1.96%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '242.65'
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('E6').Value = '  -2.30%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.03'
$ws.Range('E8').Value = '  -3.48%  '
$ws.Range('E9').Value = '  -4.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0849'
$ws.Range('E10').Value = '  -2.36%  '
$ws.Range('E11').Value = '  -1.54%  '
$ws.Range('D12').Value = '2.223.61'
$ws.Range('E12').Value = '  -1.94%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.29'
$ws.Range('E13').Value = '  -5.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.808'
$ws.Range('E14').Value = '  -5.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '13.42'
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.14'
$ws.Range('E16').Value = '  -5.87%  '
$ws.Range('D17').Value = '1.935.08'
$ws.Range('E17').Value = '  -2.54%  '
$ws.Range('D18').Value = '36.394.66'
$ws.Range('E18').Value = '  +0.32%  '
$ws.Range('B19').Value = 'Litecoin'
$ws.Range('C19').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.12'
$ws.Range('E19').Value = '  -1.92%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0862'
$ws.Range('E20').Value = '  -4.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '227.34'
$ws.Range('E21').Value = '  -3.02%  '
$ws.Range('E22').Value = '  -5.71%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('E24').Value = '  -6.74%  '
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.17'
$ws.Range('E26').Value = '  -6.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.95'
$ws.Range('E27').Value = '  -2.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.133'
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.22'
$ws.Range('E29').Value = '  -3.31%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.117'
$ws.Range('E30').Value = '  -2.27%  '
$ws.Range('E31').Value = '  -7.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.54'
$ws.Range('E32').Value = '  -6.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0617'
$ws.Range('E33').Value = '  -5.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.15'
$ws.Range('E34').Value = '  -6.28%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.08'
$ws.Range('E36').Value = '  -0.27%  '
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('E38').Value = '  -1.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.12'
$ws.Range('E39').Value = '  +6.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0989'
$ws.Range('E40').Value = '  +2.82%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.89'
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('E42').Value = '  -2.04%  '
$ws.Range('E43').Value = '  -5.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.59'
$ws.Range('E44').Value = '  -3.97%  '
$ws.Range('D45').Value = '1.337.77'
$ws.Range('E45').Value = '  -2.06%  '
$ws.Range('E46').Value = '  -6.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.34'
$ws.Range('E47').Value = '  -5.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.12'
$ws.Range('E48').Value = '  -4.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.82'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('D50').Value = '2.114.58'
$ws.Range('E50').Value = '  -1.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.12'
$ws.Range('E51').Value = '  -5.02%  '
